$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the "conteo" helper table (D:F, rows 2-6) three columns to the
#    right (G:I), freeing up D:F and making room for a narrow spacer
#    column C. Column-by-column copy preserves each cell's existing style.
$rows = 2,3,4,5,6
foreach ($r in $rows) {
  $ws.Range("D$r").Copy($ws.Range("G$r"))
  $ws.Range("E$r").Copy($ws.Range("H$r"))
  $ws.Range("F$r").Copy($ws.Range("I$r"))
}

# The COUNTIF formulas referenced the old D-column lookup cell; point them
# at the new G-column location now that the lookup values moved there too.
$ws.Range("I2").Formula = "=COUNTIF(B1:B32,G2)"
$ws.Range("I3").Formula = "=COUNTIF(B2:B33,G3)"
$ws.Range("I4").Formula = "=COUNTIF(B3:B34,G4)"
$ws.Range("I5").Formula = "=COUNTIF(B4:B35,G5)"
$ws.Range("I6").Formula = "=COUNTIF(B5:B36,G6)"

# Remove the now-vacated D:F cells entirely (contents + formatting) so
# they don't leave stray styled-but-empty cells behind.
$ws.Range("D2:F6").Clear()

# 2. Insert a narrow spacer column at C (width ~4) between the document
#    list (A:B) and the relocated helper table (G:I).
$ws.Columns("C").ColumnWidth = 3.2

# 3. Fix a couple of rows whose document-name / responsible cells were
#    still using the old plain style instead of the highlighted one used
#    throughout the rest of the list.
$fixCells = "A15","B15","A20","B20"
foreach ($c in $fixCells) {
  $ws.Range($c).Font.Name = "Cambria"
  $ws.Range($c).Font.Size = 11
  $ws.Range($c).Interior.Color = 65535
}

# 4. Update the active selection to match the refreshed sheet.
[void]$ws.Range("E9").Select()
